$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8549935817718506
$ws.Range("B1").Value = 2.867308616638184
$ws.Range("C1").Value = 8.781153678894043
$ws.Range("D1").Value = 2.026142358779907
$ws.Range("E1").Value = 1.147864580154419
